# Kandidat_2022_vars.xlsx edit:
# - Repurpose the "Jobber i Oslo" row (row 8) of the variable-reference
#   table on Sheet1 into a new factor variable describing where the
#   respondent's workplace is located (Oslo / Viken / other).
# - Update the active selection to the changed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 columns: A=Spørsmålstekst, B=Variabel, C=Svartype, D=Arkfanetittel, E=Kommentar
$ws.Range("A8").Value = "Hvor er arbeidsstedet ditt?"
$ws.Range("B8").Value = "oslo_viken_annet"
$ws.Range("C8").Value = "fordeling"
$ws.Range("D8").Value = "Arbeidssted"

# Match the saved selection state: A8 active cell, A8:D8 selected.
$ws.Activate()
$ws.Range("A8:D8").Select()
